$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "Marco sala"
$ws.Range("B21").Value = "Stefano Tita | Clitoriders"
$ws.Range("C21").Value = "Luca Frasca | Clitoriders"
$ws.Range("D21").Value = "Giacomo Gasparini | MAI UNA GIOIA"
$ws.Range("E21").Value = "daniel pedrotti | iMontagna"
$ws.Range("F21").Value = "Daniele Dalbosco | SdrumALA"
